# Generate Report for Handback
# The a33124e5-5ffc-40cc-b6c3-777194ec2b40.md file has now been successfully
# handed back (in sync with en-US) for both the zh-cn and de-de locales.
# Update the localization-status report accordingly:
#   - Overview sheet: zh-cn / de-de status columns for that file
#   - zh-cn sheet: Status, Latest Handback DateTime, Error Detail for that file
#   - de-de sheet: Status, Latest Handback DateTime, Error Detail for that file

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- "Overview" sheet (row 3 = a33124e5-5ffc-40cc-b6c3-777194ec2b40.md) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack

# --- "zh-cn" sheet (row 3 = a33124e5-5ffc-40cc-b6c3-777194ec2b40.md) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $handedBack
$wsZhCn.Range("K3").Value = "2016-09-07 07:08:56"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- "de-de" sheet (row 3 = a33124e5-5ffc-40cc-b6c3-777194ec2b40.md) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $handedBack
$wsDeDe.Range("K3").Value = "2016-09-07 07:09:12"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839
